# Updated cryptos list - apply scraped price/volume changes to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.084.69'
$ws.Range('E2').Value = '  -0.25%  '
$ws.Range('D3').Value = '3.458.70'
$ws.Range('E3').Value = '  -0.78%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '579.11'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -1.12%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '149.43'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +1.42%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('E8').Value = '  +0.10%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '7.87'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +2.49%  '
$ws.Range('E10').Value = '  -1.91%  '
$ws.Range('E11').Value = '  +2.33%  '
$ws.Range('D12').Value = '4.051.05'
$ws.Range('E13').Value = '  +2.34%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '28.57'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -4.00%  '
$ws.Range('D15').Value = '3.462.87'
$ws.Range('E15').Value = '  -0.97%  '
$ws.Range('E16').Value = '  -1.29%  '
$ws.Range('D17').Value = '63.138.21'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '6.45'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +2.67%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '14.47'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +1.13%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '9.18'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -2.76%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '388.58'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -1.08%  '
$ws.Range('E22').Value = '  -0.54%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '74.67'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -0.80%  '
$ws.Range('E24').Value = '  +0.06%  '
$ws.Range('D25').Value = '3.594.71'
$ws.Range('E25').Value = '  -0.98%  '
$ws.Range('E26').Value = '  -3.25%  '
$ws.Range('E27').Value = '  -1.66%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.68'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -2.46%  '
$ws.Range('E29').Value = '  +0.16%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.07'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -2.15%  '
$ws.Range('E31').Value = '  -1.84%  '
$ws.Range('E32').Value = '  +0.02%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '23.38'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -2.00%  '
$ws.Range('E34').Value = '  -6.24%  '
$ws.Range('E35').Value = '  +4.00%  '
$ws.Range('E36').Value = '  +0.01%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '32.04'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -1.78%  '
$ws.Range('E38').Value = '  -1.74%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '170.42'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -0.72%  '
$ws.Range('D40').Value = '3.494.80'
$ws.Range('E40').Value = '  -0.85%  '
$ws.Range('E41').Value = '  +0.91%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.794'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -0.98%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '42.86'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +0.89%  '
$ws.Range('E44').Value = '  -1.70%  '
$ws.Range('E45').Value = '  -1.71%  '
$ws.Range('E46').Value = '  -2.90%  '
$ws.Range('D47').Value = '2.584.85'
$ws.Range('E47').Value = '  -1.05%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.32'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +0.54%  '
$ws.Range('E49').Value = '  +2.00%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '22.66'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -4.70%  '
$ws.Range('E51').Value = '  -0.01%  '
